# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 764
$ws1.Range("F3").Value = 23
$ws1.Range("F5").Value = 35
$ws1.Range("F7").Value = 3626
$ws1.Range("F9").Value = 4246
$ws1.Range("F10").Value = 485
$ws1.Range("F11").Value = 1057

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 764
$ws4.Range("F3").Value = 23
$ws4.Range("F5").Value = 35
$ws4.Range("F8").Value = 3626
$ws4.Range("F10").Value = 4246
$ws4.Range("F11").Value = 485
$ws4.Range("F12").Value = 1057
